$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 16:22"

# --- Update Estados Unidos (row 4) ---
$ws.Range("B4").Value = 740151
$ws.Range("C4").Value = 1359
$ws.Range("D4").Value = 68456
$ws.Range("E4").Value = 632627
$ws.Range("G4").Value = 54
$ws.Range("H4").Value = 39068

# --- Update Alemania (row 8) ---
$ws.Range("B8").Value = 144348
$ws.Range("C8").Value = 624
$ws.Range("E8").Value = 51801
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 4547

# --- Update Suiza (row 18) ---
$ws.Range("E18").Value = 9259
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 1381

# --- Update Noruega (row 36) ---
$ws.Range("B36").Value = 7078
$ws.Range("C36").Value = 42
$ws.Range("E36").Value = 6881

# --- Update Oman (row 77) ---
$ws.Range("E77").Value = 1026
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 7

# --- Kenia overtakes Sri Lanka and Vietnam in ranking (rows 116-118) ---
# Row 116 becomes Kenia with updated figures, Sri Lanka and Vietnam shift down one
# row each, keeping their previous figures.
$ws.Range("A116").Value = "Kenia"
$ws.Range("B116").Value = 270
$ws.Range("C116").Value = 8
$ws.Range("D116").Value = 67
$ws.Range("E116").Value = 189
$ws.Range("F116").Value = 2
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 14

$ws.Range("A117").Value = "Sri Lanka"
$ws.Range("B117").Value = 269
$ws.Range("C117").Value = 15
$ws.Range("D117").Value = 96
$ws.Range("E117").Value = 166
$ws.Range("F117").Value = 1
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 7

$ws.Range("A118").Value = "Vietnam"
$ws.Range("B118").Value = 268
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 203
$ws.Range("E118").Value = 65
$ws.Range("F118").Value = 8
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 0
